# Adds the new "operation pump" UI texts to the Translation sheet
# (rows 63-69), filling in TEXT ID, TYPOGRAPHY NAME, ALIGNMENT,
# DIRECTION and GB (text) columns, matching the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$rows = @(
    @("SingleUseId68", "Typographies_button", "Center", "LTR", "DISPACHE"),
    @("SingleUseId69", "Typographies_button", "Center", "LTR", "STOP"),
    @("SingleUseId70", "Typographies_button", "Center", "LTR", "CANCEL SALE"),
    @("SingleUseId71", "Typographies_button", "Center", "LTR", "PAY SALE"),
    @("SingleUseId72", "Typography_label", "Center", "LTR", "Liters filled: <value>"),
    @("SingleUseId73", "Typography_label", "Left", "LTR", "0"),
    @("SingleUseId74", "Typography_label", "Center", "LTR", "Operation: <value>")
)

# Column F values that are purely numeric-looking text (e.g. "0") need to be
# forced to text so they keep their shared-string type instead of becoming a
# number.
$textColumnValues = @("0")

$startRow = 63
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $data[0]
    $ws.Cells.Item($r, 3).Value = $data[1]
    $ws.Cells.Item($r, 4).Value = $data[2]
    $ws.Cells.Item($r, 5).Value = $data[3]

    $fCell = $ws.Cells.Item($r, 6)
    if ($textColumnValues -contains $data[4]) {
        $fCell.NumberFormat = "@"
    }
    $fCell.Value = $data[4]
}

$wb.Save()
